$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend the header row (row 1) with two new columns: I ("I0") and J ("IF").
# Copy the existing header formatting (bold, border, centered) from H1 so the
# new header cells match the rest of the header row.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Fill in the new data column values for rows 2 and 3.
$ws.Range("I2").Value = 8
$ws.Range("J2").Value = 8
$ws.Range("I3").Value = 9
$ws.Range("J3").Value = 9
